# Add a new "localdb" command-type column/category to the hidden
# '#system' lookup sheet, used to populate Nexial's command-type /
# command-name dropdowns (data validation lists driven by named ranges).
#
# Effect:
#   1. A new column is inserted immediately before column N ("macro"),
#      shifting columns N:AC to O:AD. The new column N is populated with
#      the header "localdb" and its six function names.
#   2. A new row is inserted into column A only (the master "target" list)
#      immediately before the existing "macro" entry (row 14), shifting
#      A14:A29 down to A15:A30, keeping the list alphabetically ordered.
#      The new A14 cell is set to "localdb".
#   3. All named ranges whose column references fall at/after N are
#      shifted one column to the right, "target"'s row count grows by one,
#      and a brand-new named range "localdb" is added.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("#system")

# 1) Insert the new column before N, shifting N:AC -> O:AD.
$ws.Range("N1").EntireColumn.Insert()

$ws.Range("N1").Value = "localdb"
$ws.Range("N2").Value = "cloneTable(var,source,target)"
$ws.Range("N3").Value = "dropTables(var,tables)"
$ws.Range("N4").Value = "exportCSV(sql,output)"
$ws.Range("N5").Value = "importRecords(var,sourceDb,sql,table)"
$ws.Range("N6").Value = "purge(var)"
$ws.Range("N7").Value = "runSQLs(var,sqls)"

# 2) Insert a new entry into the "target" list (column A only) right
#    before "macro" (alphabetical order), shifting A14:A29 -> A15:A30.
$ws.Range("A14").Insert(-4121)  # -4121 = xlShiftDown
$ws.Range("A14").Value = "localdb"

# 3) Fix up the named ranges that pointed at columns N..AC (now O..AD),
#    and widen "target" by the extra row, then register "localdb".
$names = $wb.Names
foreach ($n in $names) {
    switch ($n.Name) {
        "macro"     { $n.RefersTo = "='#system'!`$O`$2:`$O`$4" }
        "mail"      { $n.RefersTo = "='#system'!`$P`$2:`$P`$2" }
        "number"    { $n.RefersTo = "='#system'!`$Q`$2:`$Q`$16" }
        "pdf"       { $n.RefersTo = "='#system'!`$R`$2:`$R`$16" }
        "rdbms"     { $n.RefersTo = "='#system'!`$S`$2:`$S`$7" }
        "redis"     { $n.RefersTo = "='#system'!`$T`$2:`$T`$10" }
        "sms"       { $n.RefersTo = "='#system'!`$U`$2:`$U`$2" }
        "sound"     { $n.RefersTo = "='#system'!`$V`$2:`$V`$5" }
        "ssh"       { $n.RefersTo = "='#system'!`$W`$2:`$W`$9" }
        "step"      { $n.RefersTo = "='#system'!`$X`$2:`$X`$4" }
        "target"    { $n.RefersTo = "='#system'!`$A`$2:`$A`$30" }
        "web"       { $n.RefersTo = "='#system'!`$Y`$2:`$Y`$127" }
        "webalert"  { $n.RefersTo = "='#system'!`$Z`$2:`$Z`$8" }
        "webcookie" { $n.RefersTo = "='#system'!`$AA`$2:`$AA`$8" }
        "ws"        { $n.RefersTo = "='#system'!`$AB`$2:`$AB`$17" }
        "ws.async"  { $n.RefersTo = "='#system'!`$AC`$2:`$AC`$8" }
        "xml"       { $n.RefersTo = "='#system'!`$AD`$2:`$AD`$21" }
    }
}

$wb.Names.Add("localdb", "='#system'!`$N`$2:`$N`$7")
